$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.3674614980063495
$ws.Range("C2").Value = 0.1346372942852554
$ws.Range("D2").Value = 0.06677200631238023
$ws.Range("E2").Value = 0.4060931560985779
$ws.Range("F2").Value = 1.862171915237411
$ws.Range("I2").Value = 1.31576150524765
$ws.Range("K2").Value = 0.4976014749058493
$ws.Range("B3").Value = 0.3359955865642803
$ws.Range("C3").Value = 0.1206964374742654
$ws.Range("D3").Value = 0.06609435825886223
$ws.Range("E3").Value = 0.3542110107287328
$ws.Range("F3").Value = 1.797652293065028
$ws.Range("I3").Value = 1.279320014550265
$ws.Range("K3").Value = 0.4521574689580348
$ws.Range("B4").Value = 0.3169488502684885
$ws.Range("C4").Value = 0.1122080978620374
$ws.Range("D4").Value = 0.0656723034989426
$ws.Range("E4").Value = 0.3224739158650607
$ws.Range("F4").Value = 1.758753930341726
$ws.Range("I4").Value = 1.25735127456737
$ws.Range("K4").Value = 0.4245923913640866
$ws.Range("B5").Value = 0.3092554506754652
$ws.Range("C5").Value = 0.1087665882742783
$ws.Range("D5").Value = 0.06549876890305839
$ws.Range("E5").Value = 0.3095682165934619
$ws.Range("F5").Value = 1.74308102942878
$ws.Range("I5").Value = 1.24849993428306
$ws.Range("K5").Value = 0.413443442269056
$ws.Range("B6").Value = 0.3079820794116301
$ws.Range("C6").Value = 0.1081961771713225
$ws.Range("D6").Value = 0.06546985916806847
$ws.Range("E6").Value = 0.3074268207167989
$ws.Range("F6").Value = 1.740489285127609
$ws.Range("I6").Value = 1.247036249712409
$ws.Range("K6").Value = 0.4115972184584109
$ws.Range("B7").Value = 0.3168448184519264
$ws.Range("C7").Value = 0.1121616138887305
$ws.Range("D7").Value = 0.06566996945688075
$ws.Range("E7").Value = 0.3222997572284783
$ws.Range("F7").Value = 1.758541839760454
$ws.Range("I7").Value = 1.257231494432872
$ws.Range("K7").Value = 0.42444169359851
$ws.Range("B8").Value = 0.3565550016707562
$ws.Range("C8").Value = 0.1298154225422081
$ws.Range("D8").Value = 0.06653956873041267
$ws.Range("E8").Value = 0.3881780573897657
$ws.Range("F8").Value = 1.839775648506247
$ws.Range("I8").Value = 1.303111429933892
$ws.Range("K8").Value = 0.4818617372618235
$ws.Range("B9").Value = 0.4366222356236165
$ws.Range("C9").Value = 0.1650208782557456
$ws.Range("D9").Value = 0.06819933090917374
$ws.Range("E9").Value = 0.518423950626854
$ws.Range("F9").Value = 2.004853372606021
$ws.Range("I9").Value = 1.396362469679786
$ws.Range("K9").Value = 0.5971885363218519
$ws.Range("B10").Value = 0.4968292398425262
$ws.Range("C10").Value = 0.1912740684736889
$ws.Range("D10").Value = 0.06939370025246205
$ws.Range("E10").Value = 0.6149319374781896
$ws.Range("F10").Value = 2.129796567972051
$ws.Range("I10").Value = 1.46695804333099
$ws.Range("K10").Value = 0.6836555483528173
$ws.Range("B11").Value = 0.5245289083201499
$ws.Range("C11").Value = 0.203308303327475
$ws.Range("D11").Value = 0.06993219330988865
$ws.Range("E11").Value = 0.6590511050090555
$ws.Range("F11").Value = 2.187462106484929
$ws.Range("I11").Value = 1.499545009991365
$ws.Range("K11").Value = 0.7233856272114849
$ws.Range("B12").Value = 0.5350634907246103
$ws.Range("C12").Value = 0.2078790215188349
$ws.Range("D12").Value = 0.07013545783453168
$ws.Range("E12").Value = 0.6757920739398315
$ws.Range("F12").Value = 2.209419891174662
$ws.Range("I12").Value = 1.511954190808567
$ws.Range("K12").Value = 0.738488502314965
$ws.Range("B13").Value = 0.5327926550976088
$ws.Range("C13").Value = 0.2068940227615315
$ws.Range("D13").Value = 0.07009170944251508
$ws.Range("E13").Value = 0.6721850501930078
$ws.Range("F13").Value = 2.20468547152359
$ws.Range("I13").Value = 1.509278554078136
$ws.Range("K13").Value = 0.7352332334718881
$ws.Range("B14").Value = 0.5253946831734027
$ws.Range("C14").Value = 0.2036840635822728
$ws.Range("D14").Value = 0.06994892880323533
$ws.Range("E14").Value = 0.6604276972341836
$ws.Range("F14").Value = 2.189266148174369
$ws.Range("I14").Value = 1.500564526718023
$ws.Range("K14").Value = 0.7246269822525164
$ws.Range("B15").Value = 0.5208691308149298
$ws.Range("C15").Value = 0.2017196575117737
$ws.Range("D15").Value = 0.0698613880445933
$ws.Range("E15").Value = 0.6532304893943319
$ws.Range("F15").Value = 2.179837197392146
$ws.Range("I15").Value = 1.495235982498386
$ws.Range("K15").Value = 0.7181379260844665
$ws.Range("B16").Value = 0.4950253117589796
$ws.Range("C16").Value = 0.1904894832386503
$ws.Range("D16").Value = 0.06935841496872541
$ws.Range("E16").Value = 0.6120532520446744
$ws.Range("F16").Value = 2.126044825968137
$ws.Range("I16").Value = 1.46483802511645
$ws.Range("K16").Value = 0.6810671472080116
$ws.Range("B17").Value = 0.4792510712150886
$ws.Range("C17").Value = 0.183623897271957
$ws.Range("D17").Value = 0.06904865079469147
$ws.Range("E17").Value = 0.5868499337546922
$ws.Range("F17").Value = 2.093258495743044
$ws.Range("I17").Value = 1.446311801886054
$ws.Range("K17").Value = 0.6584275321428379
$ws.Range("B18").Value = 0.4702074018721021
$ws.Range("C18").Value = 0.1796835811859978
$ws.Range("D18").Value = 0.06887002524302943
$ws.Range("E18").Value = 0.5723739275612019
$ws.Range("F18").Value = 2.074478504083658
$ws.Range("I18").Value = 1.435700425350177
$ws.Range("K18").Value = 0.6454430153353314
$ws.Range("B19").Value = 0.4671503786114215
$ws.Range("C19").Value = 0.1783509190740915
$ws.Range("D19").Value = 0.06880946581576097
$ws.Range("E19").Value = 0.567475999510421
$ws.Range("F19").Value = 2.068133229986813
$ws.Range("I19").Value = 1.432115186909073
$ws.Range("K19").Value = 0.6410530373226493
$ws.Range("B20").Value = 0.4809272340359882
$ws.Range("C20").Value = 0.184353858595415
$ws.Range("D20").Value = 0.06908167274365695
$ws.Range("E20").Value = 0.5895307515949071
$ws.Range("F20").Value = 2.09674058624995
$ws.Range("I20").Value = 1.448279344114781
$ws.Range("K20").Value = 0.6608337012244192
$ws.Range("B21").Value = 0.527566412292316
$ws.Range("C21").Value = 0.204626533255805
$ws.Range("D21").Value = 0.06999088425924072
$ws.Range("E21").Value = 0.6638801718571159
$ws.Range("F21").Value = 2.193791874659922
$ws.Range("I21").Value = 1.503122159732584
$ws.Range("K21").Value = 0.7277407139965533
$ws.Range("B22").Value = 0.5583122419028257
$ws.Range("C22").Value = 0.2179554705475368
$ws.Range("D22").Value = 0.07058133181703852
$ws.Range("E22").Value = 0.7126712962469952
$ws.Range("F22").Value = 2.257927289824607
$ws.Range("I22").Value = 1.539369063414668
$ws.Range("K22").Value = 0.7718066299608211
$ws.Range("B23").Value = 0.541878226248059
$ws.Range("C23").Value = 0.2108341331749841
$ws.Range("D23").Value = 0.07026652974285952
$ws.Range("E23").Value = 0.6866113771413183
$ws.Range("F23").Value = 2.223631675065661
$ws.Range("I23").Value = 1.519986036632417
$ws.Range("K23").Value = 0.7482565151435665
$ws.Range("B24").Value = 0.4801693625255439
$ws.Range("C24").Value = 0.1840238219356536
$ws.Range("D24").Value = 0.06906674519888867
$ws.Range("E24").Value = 0.5883187113625326
$ws.Range("F24").Value = 2.095166117392637
$ws.Range("I24").Value = 1.447389695021684
$ws.Range("K24").Value = 0.659745774772631
$ws.Range("B25").Value = 0.4147224710028468
$ws.Range("C25").Value = 0.1554309907685081
$ws.Range("D25").Value = 0.06775491684414447
$ws.Range("E25").Value = 0.4830581212895595
$ws.Range("F25").Value = 1.95956267409062
$ws.Range("I25").Value = 1.370775727204375
$ws.Range("K25").Value = 0.56569010782502

Write-Output "Applied 168 cell updates"
